# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the header style already used by row 1 (A1:AC1) by copying its
# formatting onto the three new header cells before writing their text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, bordered, centered - style
# index 1 in the original file) onto the new header cells, then set text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record (77-85-0) repeated on every player row, 2 through 58.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 85   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
